$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.009.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.059.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.50'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +13.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '60.55'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  +1.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0802'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.107'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.358.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('E15').Value = '  -2.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.060.78'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.979.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0935'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +13.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('E21').Value = '  +7.35%  '
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.50%  '
$ws.Range('E28').Value = '  -4.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('E34').Value = '  +6.74%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0877'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.92%  '
$ws.Range('E37').Value = '  -5.51%  '
$ws.Range('E38').Value = '  -3.71%  '
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  +23.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.06'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.25%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.10%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +50.01%  '
$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.37'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +11.70%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.35%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.300.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.28%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.92'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.50%  '
$ws.Range('B51').Value = 'Gas'
$ws.Range('C51').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '13.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -49.24%  '
